# Readmes added to templates
# Adds a new "README" sheet in front of the existing "Classification" sheet,
# matching the GBIF IPT checklist template layout.

$wb = $excel.ActiveWorkbook
$classification = $wb.Worksheets.Item("Classification")

# Insert the new sheet before Classification -> it becomes sheet index 1,
# Classification becomes index 2 (mirrors workbook.xml sheetId/r:id layout).
$ws = $wb.Worksheets.Add($classification)
$ws.Name = "README"

# ---- column widths -------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.14
$ws.Columns.Item(2).ColumnWidth = 84.5

# ---- cell text -------------------------------------------------------
$ws.Range("B1").Value = "GBIF IPT Template:"
$ws.Range("B2").Value = "Checklist Data"
$ws.Range("B4").Value = "Use this template for filling in a species list or taxonomy. Upload the template to the IPT where it can be published in Darwin Core Archive (DwC-A) format. Note this template must be mapped to the Taxon Core in the IPT."

$ws.Range("A6").Value = "Sheet Name"
$ws.Range("B6").Value = "Purpose"

$ws.Range("A7").Value = "Classification"
$ws.Range("B7").Value = "This sheet is used to record the species list. For filling in a taxonomy, make sure the classification uses a parent ID (parentNameUsageID) related to the primary taxonID as the means of recording the taxonomic hierarchy."

$ws.Range("B9").Value = "Notes:"

$ws.Range("B10").Value = "#1. The header row shows required and recommended terms. Hover over the cell to find out if it's required or recommended, and to obtain a definition of the term examples. "
$ws.Range("B11").Value = "#2. Additional columns can be added, but you should use DwC term names: http://rs.tdwg.org/dwc/terms/"
$ws.Range("B12").Value = "#3. Columns can be reordered, but the header name (equal to a DwC term name) cannot be changed."

# Bold the "#1." / "#2." / "#3." lead-ins (rich text runs).
$ws.Range("B10").Characters(1, 3).Font.Bold = $true
$ws.Range("B11").Characters(1, 3).Font.Bold = $true
$ws.Range("B12").Characters(1, 3).Font.Bold = $true

# ---- row heights -------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 20
$ws.Rows.Item(4).RowHeight = 45
$ws.Rows.Item(6).RowHeight = 16
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30

# ---- header/title box fill + borders (rows 1-6, cols A:B) --------------
# (OOXML indexed-color-table 43 == "pale yellow" #FFFF99, which this engine's
# ColorIndex palette exposes at offset 36.)
$headerBox = $ws.Range("A1:B6")
$headerBox.Interior.ColorIndex = 36

$ws.Range("A1:A5").Borders.Item(7).LineStyle = 1
$ws.Range("B1:B2").Borders.Item(7).LineStyle = 1
$ws.Range("B3:B5").Borders.Item(10).LineStyle = 1
$ws.Range("A1:B1").Borders.Item(8).LineStyle = 1
$ws.Range("B2").Borders.Item(9).LineStyle = 1

$ws.Range("A6").Borders.Item(7).LineStyle = 1
$ws.Range("A6").Borders.Item(9).LineStyle = 1
$ws.Range("A6").Borders.Item(9).Weight = -4138
$ws.Range("B6").Borders.Item(10).LineStyle = 1
$ws.Range("B6").Borders.Item(9).LineStyle = 1
$ws.Range("B6").Borders.Item(9).Weight = -4138

$ws.Range("A8").Borders.Item(7).LineStyle = 1
$ws.Range("A9").Borders.Item(7).LineStyle = 1

# Materialise the otherwise-empty cells in column A (rows 10-13) and B13 so
# they appear in the sheet XML the same way the source workbook has them.
$ws.Range("A10:A13").Font.Bold = $false
$ws.Range("B13").Font.Bold = $false

# ---- fonts ---------------------------------------------------------
$ws.Range("B1").Font.Bold = $true

$ws.Range("B2").Font.Bold = $true
$ws.Range("B2").Font.Size = 16

$ws.Range("B3").Font.Bold = $true
$ws.Range("B3").Font.Size = 10
$ws.Range("B3").Font.ColorIndex = 16
$ws.Range("B3").HorizontalAlignment = -4152

$ws.Range("B4").Font.Italic = $true
$ws.Range("B4:B5").WrapText = $true

$ws.Range("B5").Font.Italic = $true
$ws.Range("B5").Font.Size = 11

$ws.Range("A6").Font.Bold = $true
$ws.Range("A6").HorizontalAlignment = -4131
$ws.Range("B6").Font.Bold = $true
$ws.Range("B6").HorizontalAlignment = -4131

$ws.Range("A7").Font.Bold = $true
$ws.Range("A7").VerticalAlignment = -4160
$ws.Range("A8").Font.Bold = $true
$ws.Range("A8").VerticalAlignment = -4160

$ws.Range("B7:B8").Font.Size = 12
$ws.Range("B7:B8").WrapText = $true

$ws.Range("A9").Font.Size = 12

$ws.Range("B9").Font.Bold = $true
$ws.Range("B9").WrapText = $true

$ws.Range("B10:B12").WrapText = $true
$ws.Range("B11:B12").Font.Size = 12

# ---- final selection / active cell --------------------------------
$ws.Range("B14").Select()

# Window geometry (cosmetic - best effort)
$wb.Windows.Item(1).WindowState = -4143
